$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (property column) updates
$ws.Range("B2").Value = "iaest-measure:situacion-preferente"
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("F2").Value = "sdmx-dimension:refArea"
$ws.Range("H2").Value = "iaest-measure:lugar-trabajo-o-estudio"

# Row 3 (type column) updates
$ws.Range("B3").Value = "medida"
$ws.Range("D3").Value = "dim"
$ws.Range("H3").Value = "medida"

# Row 4 (datatype / URI column) updates
$ws.Range("B4").Value = "xsd:int"
$ws.Range("D4").Value = "URI-Municipio"
$ws.Range("F4").Value = "URI-Comunidad"
$ws.Range("H4").Value = "xsd:int"

# Row 5 (mapping file column) removed entirely
$ws.Range("A5:J5").EntireRow.Delete()
